$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF), styled like the existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the neighboring header cell (H1) onto the new
# header cells so they pick up the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-29
$data = @{
    2  = @(8, 8)
    3  = @(4, 4)
    4  = @(7, 7)
    5  = @(5, 5)
    6  = @(7, 7)
    7  = @(7, 8)
    8  = @(8, 8)
    9  = @(16, 16)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(5, 5)
    16 = @(9, 9)
    17 = @(6, 6)
    18 = @(7, 7)
    19 = @(9, 9)
    20 = @(7, 7)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(7, 7)
    24 = @(8, 8)
    25 = @(5, 5)
    26 = @(7, 7)
    27 = @(6, 6)
    28 = @(6, 6)
    29 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
